$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 524 (pushes existing rows 524-555 down to 525-556)
$ws.Rows.Item(524).Insert()

# Populate the newly inserted row with the new weekly price observation
$ws.Cells.Item(524,1).Value  = 10
$ws.Cells.Item(524,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(524,3).Value  = "La Araucanía"
$ws.Cells.Item(524,4).Value  = 45267
$ws.Cells.Item(524,5).Value  = 9
$ws.Cells.Item(524,6).Value  = 100112044
$ws.Cells.Item(524,7).Value  = "Perejil"
$ws.Cells.Item(524,8).Value  = "Sin especificar"
$ws.Cells.Item(524,9).Value  = "Primera"
$ws.Cells.Item(524,10).Value = 60
$ws.Cells.Item(524,11).Value = 6000
$ws.Cells.Item(524,12).Value = 6000
$ws.Cells.Item(524,13).Value = 6000
$ws.Cells.Item(524,14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(524,15).Value = "Provincia de Cautín"
$ws.Cells.Item(524,16).Value = 2000
$ws.Cells.Item(524,17).Value = 3
$ws.Cells.Item(524,18).Value = "Hortaliza"
